# Revert "meta data data dictionary data RDB"
#
# Target table (word/document.xml), columns: 資料表編號/資料表名稱/資料表中文名稱
#   T04 row: evaluations/評分        -> stu&pro/學生與專案
#   T07 row: autobiography/自傳      -> stusort/學生排序
#   T08 row: resume/履歷             -> autobiography/自傳
#   (new row) T09: resume/履歷       -> inserted right before the old last row
#   old last row "T09" user_triggers/使用者異動 -> renumbered to T10 (its
#   name/description text stays the same)

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Robustly overwrite a table cell's text. Some cells in this workbook are
# made up of more than one run (e.g. "T0" + "7"); a plain
# `$table.Cell(r,c).Range.Text = "..."` only overwrites the first run and
# leaves any trailing runs behind, so explicitly delete the cell's existing
# contents (everything except the trailing cell-mark) before inserting the
# replacement text.
function Set-CellText($table, $row, $col, $newText) {
    $cell = $table.Cell($row, $col)
    $r = $cell.Range
    if (($r.End - $r.Start) -gt 1) {
        $delRange = $d.Range($r.Start, $r.End - 1)
        $delRange.Delete()
    }
    $fresh = $table.Cell($row, $col).Range
    $fresh.Text = $newText
}

# Row 5 = T04: evaluations/評分 -> stu&pro/學生與專案
Set-CellText $t 5 2 "stu&pro"
Set-CellText $t 5 3 "學生與專案"

# Row 8 = T07: autobiography/自傳 -> stusort/學生排序
# (the "T07" identifier cell is made up of two runs, "T0" + "7"; rewrite it
# too so it collapses back down to a single run, matching the target.)
Set-CellText $t 8 1 "T07"
Set-CellText $t 8 2 "stusort"
Set-CellText $t 8 3 "學生排序"

# Row 9 = T08: resume/履歷 -> autobiography/自傳
# (same single-run collapse for the "T08" identifier cell.)
Set-CellText $t 9 1 "T08"
Set-CellText $t 9 2 "autobiography"
Set-CellText $t 9 3 "自傳"

# Insert a new row before the current last row (row 10, old "T09") and fill
# it in with T09 / resume / 履歷.
[void]$t.Rows.Add($t.Rows.Item(10))
Set-CellText $t 10 1 "T09"
Set-CellText $t 10 2 "resume"
Set-CellText $t 10 3 "履歷"

# The old last row (previously row 10, now row 11) was "T09" / user_triggers
# / 使用者異動; renumber its identifier to T10 (name/description unchanged).
Set-CellText $t 11 1 "T10"

Write-Output "done"
